$d = $word.ActiveDocument

$replacements = @(
    @("2025-06-10 Tuesday", "2025-06-11 Wednesday"),
    @("267×4=", "218×8="),
    @("710×5=", "685×5="),
    @("118×6=", "434×3="),
    @("935×9=", "781×8="),
    @("667×9=", "876×4="),
    @("762×2=", "875×2="),
    @("178×4=", "878×4="),
    @("690×7=", "999×7="),
    @("935×2=", "663×8="),
    @("923×4=", "234×9="),
    @("907×8=", "994×7="),
    @("689×3=", "775×6="),
    @("371×7=", "935×7="),
    @("131×8=", "970×6="),
    @("191×8=", "290×3="),
    @("854×7=", "725×6="),
    @("886×6=", "691×9="),
    @("263×5=", "741×4="),
    @("363×9=", "457×5="),
    @("623×7=", "524×5="),
    @("513×7=", "339×4="),
    @("975×6=", "502×4="),
    @("136×5=", "379×2="),
    @("512×6=", "276×7="),
    @("738×2=", "844×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
